$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block of 9 rows of data to insert right after the header row (row 1).
$newRows = @(
    @(-0.1018617823719978, 0.0229074470698833, 0.0226020142436027),
    @(-0.0087048299610614, -0.0256563406437635, -0.0074830991216003),
    @(-0.0534507073462009, -0.009010262787342, 0.001527163083665),
    @(-0.0387899428606033, 0.0113010071218013, 0.0397062413394451),
    @(-0.011148290708661, -0.08338310569524759, 0.0064140851609408),
    @(-0.0229074470698833, -0.0862847194075584, 0.0215329993516206),
    @(-0.0293215326964855, -0.0226020142436027, -0.0091629782691597),
    @(-0.0369573459029197, -0.0378736443817615, -0.0192422550171613),
    @(0.0308486949652433, -0.0488692186772823, -0.0198531206697225)
)

$insertCount = $newRows.Count

# Insert blank rows starting at row 2, shifting the existing data down.
$insertRange = $ws.Range("A2:C" + (1 + $insertCount))
$insertRange.Insert()
# Drop any formatting that Insert() may have copied down from the header row.
$insertRange.ClearFormats()

# Fill the newly inserted rows with the new data.
for ($i = 0; $i -lt $insertCount; $i++) {
    $rowNum = 2 + $i
    $ws.Cells.Item($rowNum, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $newRows[$i][2]
}

# Append one additional row of data at the very end (new row 31).
$lastRow = 31
$ws.Cells.Item($lastRow, 1).Value = -0.0236710291355848
$ws.Cells.Item($lastRow, 2).Value = 0.039248090237379
$ws.Cells.Item($lastRow, 3).Value = 0.0468839071691036
